$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) End the "As features são passadas depois de and ()" paragraph with a
#    new run containing just a period.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*As features são passadas depois de and ()*") {
        $target = $p
        break
    }
}

$rng = $target.Range
# Trim off the trailing paragraph mark so Collapse(0) lands right after ")"
$rng.SetRange($rng.Start, $rng.End - 1)
$rng.Collapse(0)
$rng.InsertAfter(".")

# ---------------------------------------------------------------------------
# 2) Insert a brand-new paragraph right after it, with the "MEDIA FEATURE"
#    explanation (mixed bold/regular runs).
# ---------------------------------------------------------------------------
$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null

$newPara = $target.Next()
$newRng = $newPara.Range

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/frag.xml" pkg:contentType="application/xml"><pkg:xmlData><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:tab/><w:t xml:space="preserve">A grande sacada dos </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">MEDIA FEATURE </w:t></w:r><w:r><w:t xml:space="preserve">é que nós podemos passar </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">configurações gerais </w:t></w:r><w:r><w:t xml:space="preserve">que sempre (faça chuva ou faça sol) funcionarão de uma determinada maneira </w:t></w:r><w:r><w:t xml:space="preserve">(por exemplo: background-color, font-famyli, color, font-weight). E podemos passar features que funcionarão em situações específicas. A grande complexidade em se usar essas configurações de featura é entender o que é uma </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">configuração específica </w:t></w:r><w:r><w:t xml:space="preserve">e o que é uma </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>configuração geral.</w:t></w:r></w:p></pkg:xmlData></pkg:part></pkg:package>'
$newRng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the stray <w:lastRenderedPageBreak/> right before the
#    '@charset "UTF-8";' paragraph.
# ---------------------------------------------------------------------------
$charsetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*@charset*UTF-8*") {
        $charsetPara = $p
        break
    }
}
$cRng = $charsetPara.Range
$cRng.Find.ClearFormatting()
$cRng.Find.Execute("@charset", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cRng.Collapse(1)
$cRng.MoveStart(1, -1) | Out-Null
$cRng.Delete() | Out-Null
